# Fill in the missing (previously 0) measurement values for Lessons A and B
# so the results table matches the completed answers for Lessons A through C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01_sonora_bobcat_2018")

# --- Lesson A, specimen 1 (row 3) ---
$ws.Range("D3").Value = 1900
$ws.Range("G3").Value = 4.25
$ws.Range("J3").Value = 36.6
$ws.Range("M3").Value = -14

# --- Lesson A, specimen 2 (row 4) ---
$ws.Range("D4").Value = 1150
$ws.Range("G4").Value = 4
$ws.Range("J4").Value = 21.3
$ws.Range("M4").Value = 41.8

# --- Lesson A, specimen 3 (row 5) ---
$ws.Range("D5").Value = 875
$ws.Range("G5").Value = 5.5
$ws.Range("J5").Value = 4.4000000000000004
$ws.Range("M5").Value = 6.6

# --- Lesson B, specimen 1 (row 6) ---
$ws.Range("D6").Value = 700
$ws.Range("F6").Value = 4.5
$ws.Range("G6").Value = 4.75
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 25.3
$ws.Range("J6").Value = 36
$ws.Range("K6").Value = 40
$ws.Range("M6").Value = -41.75

# --- Lesson B, specimen 2 (row 7) ---
$ws.Range("D7").Value = 2100
$ws.Range("G7").Value = 5.5
$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 18
$ws.Range("K7").Value = 25
$ws.Range("M7").Value = 12.3
$ws.Range("N7").Value = 25.5

# --- Lesson B, specimen 3 (row 8) ---
$ws.Range("D8").Value = 800
$ws.Range("G8").Value = 4.25
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 26.3
$ws.Range("K8").Value = 35
$ws.Range("M8").Value = -4.95
$ws.Range("N8").Value = 6.55

# Update the active selection to reflect where the editor ended up (C8)
$ws.Range("C8").Select()
